$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '66.318.24'
$r.Style = 'Normal'
$r = $ws.Range('E2')
$r.NumberFormat = '@'
$r.Value = '  +0.61%  '
$r.Style = 'Normal'
$r = $ws.Range('D3')
$r.NumberFormat = '@'
$r.Value = '3.345.81'
$r.Style = 'Normal'
$r = $ws.Range('E3')
$r.NumberFormat = '@'
$r.Value = '  +1.22%  '
$r.Style = 'Normal'
$r = $ws.Range('D4')
$r.NumberFormat = '@'
$r.Value = '0.997'
$r.Style = 'Normal'
$r = $ws.Range('E4')
$r.NumberFormat = '@'
$r.Value = '  -0.38%  '
$r.Style = 'Normal'
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '585.27'
$r.Style = 'Normal'
$r = $ws.Range('E5')
$r.NumberFormat = '@'
$r.Value = '  +5.31%  '
$r.Style = 'Normal'
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '185.82'
$r.Style = 'Normal'
$r = $ws.Range('E6')
$r.NumberFormat = '@'
$r.Value = '  -0.92%  '
$r.Style = 'Normal'
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.Style = 'Normal'
$r = $ws.Range('E7')
$r.NumberFormat = '@'
$r.Value = '  -0.05%  '
$r.Style = 'Normal'
$r = $ws.Range('E8')
$r.NumberFormat = '@'
$r.Value = '  -1.03%  '
$r.Style = 'Normal'
$r = $ws.Range('E9')
$r.NumberFormat = '@'
$r.Value = '  -0.42%  '
$r.Style = 'Normal'
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '0.583'
$r.Style = 'Normal'
$r = $ws.Range('E10')
$r.NumberFormat = '@'
$r.Value = '  -0.32%  '
$r.Style = 'Normal'
$r = $ws.Range('D11')
$r.NumberFormat = '@'
$r.Value = '47.00'
$r.Style = 'Normal'
$r = $ws.Range('E11')
$r.NumberFormat = '@'
$r.Value = '  -0.50%  '
$r.Style = 'Normal'
$r = $ws.Range('E12')
$r.NumberFormat = '@'
$r.Value = '  +0.17%  '
$r.Style = 'Normal'
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '655.71'
$r.Style = 'Normal'
$r = $ws.Range('E13')
$r.NumberFormat = '@'
$r.Value = '  +8.58%  '
$r.Style = 'Normal'
$r = $ws.Range('E14')
$r.NumberFormat = '@'
$r.Value = '  -2.08%  '
$r.Style = 'Normal'
$r = $ws.Range('D15')
$r.NumberFormat = '@'
$r.Value = '3.632.82'
$r.Style = 'Normal'
$r = $ws.Range('E15')
$r.NumberFormat = '@'
$r.Value = '  -5.39%  '
$r.Style = 'Normal'
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '66.408.63'
$r.Style = 'Normal'
$r = $ws.Range('E16')
$r.NumberFormat = '@'
$r.Value = '  +0.68%  '
$r.Style = 'Normal'
$r = $ws.Range('E17')
$r.NumberFormat = '@'
$r.Value = '  -0.12%  '
$r.Style = 'Normal'
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '17.89'
$r.Style = 'Normal'
$r = $ws.Range('E18')
$r.NumberFormat = '@'
$r.Value = '  -0.36%  '
$r.Style = 'Normal'
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '3.339.03'
$r.Style = 'Normal'
$r = $ws.Range('E19')
$r.NumberFormat = '@'
$r.Value = '  +1.76%  '
$r.Style = 'Normal'
$r = $ws.Range('E20')
$r.NumberFormat = '@'
$r.Value = '  +0.67%  '
$r.Style = 'Normal'
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '0.899'
$r.Style = 'Normal'
$r = $ws.Range('E21')
$r.NumberFormat = '@'
$r.Value = '  -0.74%  '
$r.Style = 'Normal'
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '17.70'
$r.Style = 'Normal'
$r = $ws.Range('E22')
$r.NumberFormat = '@'
$r.Value = '  -3.46%  '
$r.Style = 'Normal'
$r = $ws.Range('E23')
$r.NumberFormat = '@'
$r.Value = '  +0.34%  '
$r.Style = 'Normal'
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '100.32'
$r.Style = 'Normal'
$r = $ws.Range('E24')
$r.NumberFormat = '@'
$r.Value = '  +0.26%  '
$r.Style = 'Normal'
$r = $ws.Range('E25')
$r.NumberFormat = '@'
$r.Value = '  +1.05%  '
$r.Style = 'Normal'
$r = $ws.Range('E26')
$r.NumberFormat = '@'
$r.Value = '  +1.20%  '
$r.Style = 'Normal'
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '9.60'
$r.Style = 'Normal'
$r = $ws.Range('E27')
$r.NumberFormat = '@'
$r.Value = '  +0.43%  '
$r.Style = 'Normal'
$r = $ws.Range('D28')
$r.NumberFormat = '@'
$r.Value = '32.08'
$r.Style = 'Normal'
$r = $ws.Range('E28')
$r.NumberFormat = '@'
$r.Value = '  +5.98%  '
$r.Style = 'Normal'
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '8.55'
$r.Style = 'Normal'
$r = $ws.Range('E29')
$r.NumberFormat = '@'
$r.Value = '  -1.65%  '
$r.Style = 'Normal'
$r = $ws.Range('D30')
$r.NumberFormat = '@'
$r.Value = '6.84'
$r.Style = 'Normal'
$r = $ws.Range('E30')
$r.NumberFormat = '@'
$r.Value = '  +1.76%  '
$r.Style = 'Normal'
$r = $ws.Range('D31')
$r.NumberFormat = '@'
$r.Value = '602.30'
$r.Style = 'Normal'
$r = $ws.Range('E31')
$r.NumberFormat = '@'
$r.Value = '  +3.92%  '
$r.Style = 'Normal'
$r = $ws.Range('D32')
$r.NumberFormat = '@'
$r.Value = '3.89'
$r.Style = 'Normal'
$r = $ws.Range('E32')
$r.NumberFormat = '@'
$r.Value = '  +0.47%  '
$r.Style = 'Normal'
$r = $ws.Range('D33')
$r.NumberFormat = '@'
$r.Value = '11.10'
$r.Style = 'Normal'
$r = $ws.Range('E33')
$r.NumberFormat = '@'
$r.Value = '  +0.26%  '
$r.Style = 'Normal'
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '3.879.82'
$r.Style = 'Normal'
$r = $ws.Range('E34')
$r.NumberFormat = '@'
$r.Value = '  +4.75%  '
$r.Style = 'Normal'
$r = $ws.Range('E35')
$r.NumberFormat = '@'
$r.Value = '  +0.91%  '
$r.Style = 'Normal'
$r = $ws.Range('E36')
$r.NumberFormat = '@'
$r.Value = '  +0.10%  '
$r.Style = 'Normal'
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '56.49'
$r.Style = 'Normal'
$r = $ws.Range('E37')
$r.NumberFormat = '@'
$r.Value = '  -0.95%  '
$r.Style = 'Normal'
$r = $ws.Range('D38')
$r.NumberFormat = '@'
$r.Value = '2.73'
$r.Style = 'Normal'
$r = $ws.Range('E38')
$r.NumberFormat = '@'
$r.Value = '  +2.48%  '
$r.Style = 'Normal'
$r = $ws.Range('E39')
$r.NumberFormat = '@'
$r.Value = '  -0.22%  '
$r.Style = 'Normal'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '33.09'
$r.Style = 'Normal'
$r = $ws.Range('E40')
$r.NumberFormat = '@'
$r.Value = '  -2.40%  '
$r.Style = 'Normal'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.0₃0702'
$r.Style = 'Normal'
$r = $ws.Range('E41')
$r.NumberFormat = '@'
$r.Value = '  -2.10%  '
$r.Style = 'Normal'
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '3.20'
$r.Style = 'Normal'
$r = $ws.Range('E42')
$r.NumberFormat = '@'
$r.Value = '  -2.38%  '
$r.Style = 'Normal'
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.341'
$r.Style = 'Normal'
$r = $ws.Range('E43')
$r.NumberFormat = '@'
$r.Value = '  +0.29%  '
$r.Style = 'Normal'
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '3.37'
$r.Style = 'Normal'
$r = $ws.Range('E44')
$r.NumberFormat = '@'
$r.Value = '  -0.22%  '
$r.Style = 'Normal'
$r = $ws.Range('E45')
$r.NumberFormat = '@'
$r.Value = '  -0.71%  '
$r.Style = 'Normal'
$r = $ws.Range('E46')
$r.NumberFormat = '@'
$r.Value = '  -1.22%  '
$r.Style = 'Normal'
$r = $ws.Range('B47')
$r.NumberFormat = '@'
$r.Value = 'ThetaToken'
$r.Style = 'Normal'
$r = $ws.Range('C47')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$r.Style = 'Normal'
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '2.56'
$r.Style = 'Normal'
$r = $ws.Range('E47')
$r.NumberFormat = '@'
$r.Value = '  -1.20%  '
$r.Style = 'Normal'
$r = $ws.Range('B48')
$r.NumberFormat = '@'
$r.Value = 'FirstDigitalUSD'
$r.Style = 'Normal'
$r = $ws.Range('C48')
$r.NumberFormat = '@'
$r.Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$r.Style = 'Normal'
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '1.00'
$r.Style = 'Normal'
$r = $ws.Range('E48')
$r.NumberFormat = '@'
$r.Value = '  +0.18%  '
$r.Style = 'Normal'
$r = $ws.Range('E49')
$r.NumberFormat = '@'
$r.Value = '  -16.04%  '
$r.Style = 'Normal'
$r = $ws.Range('E50')
$r.NumberFormat = '@'
$r.Value = '  +7.04%  '
$r.Style = 'Normal'
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '129.93'
$r.Style = 'Normal'
$r = $ws.Range('E51')
$r.NumberFormat = '@'
$r.Value = '  +5.10%  '
$r.Style = 'Normal'
